$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily refresh of the crypto price/volume table (column D = Price, column E
# = Volume(1h)); two rows also swap rank position (13/14 and 19/20), so their
# Coin name + Link columns are updated too.
#
# Column D values that look like plain decimals (e.g. "0.9980", "11.60")
# are written with a leading apostrophe, PowerShell's/Excel's classic
# "treat this as text, not a number" marker, so trailing zeros and decimal
# formatting survive exactly like the source inline strings. Values that
# already fail to parse as a number (e.g. "26.329.33", with two dots) don't
# need it.
$ws.Range("D2").Value = '26.329.33'
$ws.Range("E2").Value = '  +4.74%  '
$ws.Range("D3").Value = '1.716.81'
$ws.Range("E3").Value = '  +3.96%  '
$ws.Range("D4").Value = '''0.9980'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''240.55'
$ws.Range("E5").Value = '  +2.76%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '''0.4725'
$ws.Range("E7").Value = '  -1.07%  '
$ws.Range("D8").Value = '''0.2641'
$ws.Range("E8").Value = '  +2.89%  '
$ws.Range("D9").Value = '''0.06239'
$ws.Range("E9").Value = '  +2.17%  '
$ws.Range("D10").Value = '1.708.14'
$ws.Range("E10").Value = '  +3.99%  '
$ws.Range("D11").Value = '''0.07086'
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").Value = '''15.29'
$ws.Range("E12").Value = '  +6.08%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.5925'
$ws.Range("E13").Value = '  +2.85%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''4.429'
$ws.Range("E14").Value = '  +2.66%  '
$ws.Range("D15").Value = '''76.17'
$ws.Range("E15").Value = '  +3.37%  '
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '''0.9988'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '26.310.29'
$ws.Range("E18").Value = '  +4.70%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '''0.000006808'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '''11.60'
$ws.Range("E20").Value = '  +2.71%  '
$ws.Range("D21").Value = '1.921.54'
$ws.Range("E21").Value = '  +3.70%  '
$ws.Range("D22").Value = '''4.569'
$ws.Range("E22").Value = '  +5.75%  '
$ws.Range("D23").Value = '''8.850'
$ws.Range("E23").Value = '  +4.42%  '
$ws.Range("D24").Value = '''5.352'
$ws.Range("E24").Value = '  +2.10%  '
$ws.Range("D25").Value = '''135.62'
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").Value = '''15.21'
$ws.Range("E26").Value = '  +1.78%  '
$ws.Range("D27").Value = '''1.403'
$ws.Range("E27").Value = '  +2.09%  '
$ws.Range("D28").Value = '''1.766'
$ws.Range("E28").Value = '  +7.37%  '
$ws.Range("D29").Value = '''106.43'
$ws.Range("E29").Value = '  +2.72%  '
$ws.Range("D30").Value = '''4.043'
$ws.Range("E30").Value = '  +3.08%  '
$ws.Range("D31").Value = '''3.699'
$ws.Range("E31").Value = '  +4.51%  '
$ws.Range("D32").Value = '''0.07766'
$ws.Range("E32").Value = '  +2.09%  '
$ws.Range("D33").Value = '''0.04414'
$ws.Range("E33").Value = '  +2.81%  '
$ws.Range("D34").Value = '''2.611'
$ws.Range("E34").Value = '  +1.36%  '
$ws.Range("D35").Value = '''0.6219'
$ws.Range("E35").Value = '  +4.57%  '
$ws.Range("D36").Value = '''0.9706'
$ws.Range("E36").Value = '  +3.68%  '
$ws.Range("D37").Value = '''0.9165'
$ws.Range("E37").Value = '  +7.25%  '
$ws.Range("D38").Value = '''112.11'
$ws.Range("E38").Value = '  +13.54%  '
$ws.Range("D39").Value = '''2.412'
$ws.Range("E39").Value = '  -6.57%  '
$ws.Range("D40").Value = '''1.912'
$ws.Range("E40").Value = '  +6.72%  '
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").Value = '''0.01469'
$ws.Range("E42").Value = '  -0.93%  '
$ws.Range("D43").Value = '''0.3823'
$ws.Range("E43").Value = '  +3.70%  '
$ws.Range("D44").Value = '''5.142'
$ws.Range("E44").Value = '  +10.72%  '
$ws.Range("D45").Value = '''0.1142'
$ws.Range("E45").Value = '  +4.39%  '
$ws.Range("D46").Value = '''6.249'
$ws.Range("E46").Value = '  +2.48%  '
$ws.Range("D47").Value = '''0.05296'
$ws.Range("E47").Value = '  +1.46%  '
$ws.Range("E48").Value = '  +5.44%  '
$ws.Range("D49").Value = '''7.693'
$ws.Range("E49").Value = '  +6.48%  '
$ws.Range("D50").Value = '''1.225'
$ws.Range("E50").Value = '  +1.74%  '
$ws.Range("D51").Value = '''0.3385'
$ws.Range("E51").Value = '  +3.07%  '
